$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 50
$ws.Range("I4").Value = 60
$ws.Range("K4").Value = 60
$ws.Range("M4").Value = 54

$ws.Range("H33").Value = 540.2
$ws.Range("I33").Value = 306.64285
$ws.Range("K33").Value = 306.64285
$ws.Range("M33").Value = -77.64285000000001

$ws.Range("H62").Value = 7855
$ws.Range("I62").Value = 6572.143
$ws.Range("J62").Value = 12345
$ws.Range("K62").Value = 6572.143
$ws.Range("L62").Value = 12345
$ws.Range("M62").Value = -5948.143
$ws.Range("N62").Value = -13593

$ws.Range("H65").Value = 7855
$ws.Range("I65").Value = 6572.143
$ws.Range("J65").Value = 12345
$ws.Range("K65").Value = 32860.715
$ws.Range("L65").Value = 61725
$ws.Range("M65").Value = -29740.715
$ws.Range("N65").Value = -67965

$ws.Range("H94").Value = 800
$ws.Range("I94").Value = 800
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 800
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -349
$ws.Range("N94").ClearContents()

$ws.Range("H100").Value = 3903.276
$ws.Range("I100").Value = 3012.1875
$ws.Range("K100").Value = 3012.1875
$ws.Range("M100").Value = -2471.1875

$ws.Range("H113").Value = 8709
$ws.Range("I113").Value = 8534.6
$ws.Range("J113").Value = 8999.666999999999
$ws.Range("K113").Value = 8534.6
$ws.Range("L113").Value = 8999.666999999999
$ws.Range("M113").Value = -5280.6
$ws.Range("N113").Value = -15507.667

$ws.Range("H135").Value = 2406
$ws.Range("I135").Value = 1071.3334
$ws.Range("K135").Value = 9642.000599999999
$ws.Range("M135").Value = -7107.000599999999

$ws.Range("H137").Value = 2418034.8
$ws.Range("I137").Value = 2101.25
$ws.Range("J137").Value = 7940168
$ws.Range("K137").Value = 6303.75
$ws.Range("L137").Value = 23820504
$ws.Range("M137").Value = -3753.75
$ws.Range("N137").Value = -23825604

$ws.Range("H141").Value = 5834.273
$ws.Range("I141").Value = 4726.9116
$ws.Range("J141").Value = 9599.299999999999
$ws.Range("K141").Value = 14180.7348
$ws.Range("L141").Value = 28797.9
$ws.Range("M141").Value = -9000.734800000002
$ws.Range("N141").Value = -39157.89999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 735.6
$ws.Range("I5").Value = 719.55554
$ws.Range("K5").Value = 719.55554
$ws.Range("M5").Value = -607.55554

$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

$ws.Range("H25").Value = 1000
$ws.Range("I25").Value = 1000
$ws.Range("K25").Value = 1000
$ws.Range("M25").Value = -598

$ws.Range("H32").Value = 19539932
$ws.Range("I32").Value = 19297982
$ws.Range("K32").Value = 19297982
$ws.Range("M32").Value = -19297695

$ws.Range("H45").Value = 4869.933
$ws.Range("I45").Value = 5281.727
$ws.Range("K45").Value = 5281.727
$ws.Range("M45").Value = -4904.727

$ws.Range("H74").Value = 2797.6667
$ws.Range("I74").Value = 2765.3635
$ws.Range("K74").Value = 2765.3635
$ws.Range("M74").Value = -1891.3635

$ws.Range("H77").Value = 2797.6667
$ws.Range("I77").Value = 2765.3635
$ws.Range("K77").Value = 13826.8175
$ws.Range("M77").Value = -9458.817499999999

$ws.Range("H88").Value = 1287.1
$ws.Range("I88").Value = 763
$ws.Range("J88").Value = 1811.2
$ws.Range("K88").Value = 763
$ws.Range("L88").Value = 1811.2
$ws.Range("M88").Value = -357
$ws.Range("N88").Value = -2623.2

$ws.Range("H91").Value = 1287.1
$ws.Range("I91").Value = 763
$ws.Range("J91").Value = 1811.2
$ws.Range("K91").Value = 763
$ws.Range("L91").Value = 1811.2
$ws.Range("M91").Value = 641
$ws.Range("N91").Value = -4619.2

$ws.Range("H102").Value = 1650.8182
$ws.Range("I102").Value = 1440.9
$ws.Range("K102").Value = 1440.9
$ws.Range("M102").Value = 181.0999999999999

$ws.Range("H110").Value = 2329.8
$ws.Range("I110").Value = 2037.25
$ws.Range("K110").Value = 2037.25
$ws.Range("M110").Value = 7.75

$ws.Range("H122").Value = 4610.7646
$ws.Range("I122").Value = 3064.7144
$ws.Range("J122").Value = 5693
$ws.Range("K122").Value = 9194.143199999999
$ws.Range("L122").Value = 17079
$ws.Range("M122").Value = -6744.143199999999
$ws.Range("N122").Value = -21979

$ws.Range("H132").Value = 5801.091
$ws.Range("I132").Value = 7502.5
$ws.Range("K132").Value = 22507.5
$ws.Range("M132").Value = -19977.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 735.6
$ws.Range("I4").Value = 719.55554
$ws.Range("K4").Value = 719.55554
$ws.Range("M4").Value = -604.55554

$ws.Range("H26").Value = 82286
$ws.Range("I26").Value = 38666.332
$ws.Range("K26").Value = 38666.332
$ws.Range("M26").Value = -38374.332

$ws.Range("H86").Value = 1388.7273
$ws.Range("I86").Value = 919.55554
$ws.Range("K86").Value = 919.55554
$ws.Range("M86").Value = 203.44446

$ws.Range("H89").Value = 1388.7273
$ws.Range("I89").Value = 919.55554
$ws.Range("K89").Value = 4597.7777
$ws.Range("M89").Value = 1018.2223

$ws.Range("H99").Value = 2321.6365
$ws.Range("I99").Value = 2025.0588
$ws.Range("K99").Value = 2025.0588
$ws.Range("M99").Value = -527.0588

$ws.Range("H105").Value = 2888.3333
$ws.Range("I105").Value = 2333.3333
$ws.Range("K105").Value = 2333.3333
$ws.Range("M105").Value = -586.3332999999998

$ws.Range("H141").Value = 155749
$ws.Range("J141").Value = 155749
$ws.Range("L141").Value = 155749
$ws.Range("N141").Value = -166109

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9413.839
$ws.Range("I31").Value = 4880.25
$ws.Range("K31").Value = 4880.25
$ws.Range("M31").Value = -4585.25

$ws.Range("H34").Value = 9413.839
$ws.Range("I34").Value = 4880.25
$ws.Range("K34").Value = 4880.25
$ws.Range("M34").Value = -4678.25

$ws.Range("H140").Value = 699999
$ws.Range("J140").Value = 699999
$ws.Range("L140").Value = 699999
$ws.Range("N140").Value = -710359

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 134
$ws.Range("I26").Value = 134
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 402
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -114
$ws.Range("N26").ClearContents()

$ws.Range("H33").Value = 958.2857
$ws.Range("I33").Value = 627.75
$ws.Range("J33").Value = 1399
$ws.Range("K33").Value = 3766.5
$ws.Range("L33").Value = 8394
$ws.Range("M33").Value = -3483.5
$ws.Range("N33").Value = -8960

$ws.Range("H60").Value = 747.3333
$ws.Range("I60").Value = 947.6667
$ws.Range("K60").Value = 2843.0001
$ws.Range("M60").Value = -2592.0001

$ws.Range("H62").Value = 3937.375
$ws.Range("I62").Value = 2249.5
$ws.Range("K62").Value = 6748.5
$ws.Range("M62").Value = -6062.5

$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

$ws.Range("H65").Value = 3937.375
$ws.Range("I65").Value = 2249.5
$ws.Range("K65").Value = 20245.5
$ws.Range("M65").Value = -16813.5

$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

$ws.Range("H130").Value = 4999.8335
$ws.Range("I130").Value = 3999
$ws.Range("J130").Value = 5200
$ws.Range("K130").Value = 11997
$ws.Range("L130").Value = 15600
$ws.Range("M130").Value = -6977
$ws.Range("N130").Value = -25640

$ws.Range("H136").Value = 2192.2222
$ws.Range("I136").Value = 1966.25
$ws.Range("K136").Value = 5898.75
$ws.Range("M136").Value = -798.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2015.6923
$ws.Range("I102").Value = 1871.1666
$ws.Range("J102").Value = 3750
$ws.Range("K102").Value = 1871.1666
$ws.Range("L102").Value = 3750
$ws.Range("M102").Value = -249.1666
$ws.Range("N102").Value = -6994

$ws.Range("H122").Value = 6600.4
$ws.Range("I122").Value = 8999.5
$ws.Range("J122").Value = 5001
$ws.Range("K122").Value = 26998.5
$ws.Range("L122").Value = 15003
$ws.Range("M122").Value = -24548.5
$ws.Range("N122").Value = -19903

$ws.Range("H126").Value = 2410.95
$ws.Range("J126").Value = 2666.5
$ws.Range("L126").Value = 7999.5
$ws.Range("N126").Value = -12939.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 32
$ws.Range("I23").Value = 32
$ws.Range("K23").Value = 32
$ws.Range("M23").Value = 198

$ws.Range("H33").Value = 59999
$ws.Range("J33").Value = 59999
$ws.Range("L33").Value = 59999
$ws.Range("N33").Value = -60579

$ws.Range("H122").Value = 15581.034
$ws.Range("I122").Value = 17643.234
$ws.Range("K122").Value = 52929.702
$ws.Range("M122").Value = -50479.702

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 10338.105
$ws.Range("I96").Value = 5430.636
$ws.Range("J96").Value = 17085.875
$ws.Range("K96").Value = 5430.636
$ws.Range("L96").Value = 17085.875
$ws.Range("M96").Value = -4057.636
$ws.Range("N96").Value = -19831.875

$ws.Range("H126").Value = 2020.5294
$ws.Range("J126").Value = 3850
$ws.Range("L126").Value = 11550
$ws.Range("N126").Value = -16490
